# Change "Version 1." to "Version 2." while reproducing the exact run
# layout produced by the original author's edit:
#   - "Version" is split into two runs: "Versi" + "on"
#   - " 1." becomes " 2" (run kept before the _GoBack bookmark)
#   - a new trailing run holding "." is appended after the bookmark
#
# Word normally keeps adjacent same-formatted text in a single run, so a
# plain text replace would not reproduce the split. Instead we force a
# run boundary at the desired offset by adding a (temporary) bookmark
# there and then deleting it again - the bookmark insertion splits the
# underlying run but deleting the bookmark itself does not rejoin the
# two runs, leaving a clean split with no leftover formatting.

$d = $word.ActiveDocument

# --- Split "Version" -> "Versi" | "on" ------------------------------
# Character offset 5 is right after "Versi" or before "on".
$splitRange = $d.Range(5, 5)
$d.Bookmarks.Add("__splitmark__", $splitRange)
$d.Bookmarks("__splitmark__").Delete()

# --- Change the digit "1" to "2" -------------------------------------
# Offset 8-9 is the "1" character (" 1." starts at offset 7).
$d.Range(8, 9).Text = "2"

# --- Remove the trailing period that currently sits before the bookmark
# After the previous edit the text is "Version 2." and the period is
# at offset 9-10, immediately before the _GoBack bookmark.
$d.Range(9, 10).Delete()

# --- Re-add the period as a new run placed after the bookmark --------
$tail = $d.Range(10, 10)
$tail.InsertAfter(".")
